$wb = $excel.ActiveWorkbook

# Add new row of data (Danik's father "Manassypov") to tab2 (sheet2)
$ws2 = $wb.Worksheets.Item("tab2")
$ws2.Range("A4").Value = "Manassypov"
$ws2.Range("B4").Value = 10

# Select the new bottom cell on tab2, mirroring the prior pattern of
# selecting the cell just below the last data row
$ws2.Range("A5").Select()

# Make tab2 the active sheet/tab (was tab1 before)
$ws2.Activate()
